# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# New K values, row-by-row for data rows 2..62 (i.e. G2:G62).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(1,2,0,0,4,1,3,1,0,0,1,0,0,3,2,2,1,1,0,0,0,2,1,3,1,1,1,1,1,2,0,2,3,1,1,1,0,1,2,1,1,1,1,2,1,1,0,2,2,1,1,1,2,2,1,1,2,3,1,1,3)

$firstRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $firstRow + $i
    $ws.Range("G$row").Value = $newK[$i]
}
